# Auto-generated edit script: updates Horarios Linea 141 schedule data
# across the three worksheets (LP1912, LP1912-215, 6203-6173)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 1).Value = "Última actualización: 04:18:53"
$ws.Cells.Item(3, 1).Value = "Total filas: 13"
$ws.Cells.Item(7, 1).Value = "04:18:53"
$ws.Cells.Item(7, 2).Value = "04:45"
$ws.Cells.Item(7, 4).Value = 27
$ws.Cells.Item(8, 2).Value = "04:46"
$ws.Cells.Item(8, 3).Value = "215A_EL PATO"
$ws.Cells.Item(8, 4).Value = 54
$ws.Cells.Item(9, 1).Value = "04:18:53"
$ws.Cells.Item(9, 2).Value = "04:53"
$ws.Cells.Item(9, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(9, 4).Value = 35
$ws.Cells.Item(10, 1).Value = "04:18:53"
$ws.Cells.Item(10, 2).Value = "05:16"
$ws.Cells.Item(10, 3).Value = "17_ROMERO"
$ws.Cells.Item(10, 4).Value = 58
$ws.Cells.Item(11, 1).Value = "04:18:53"
$ws.Cells.Item(11, 2).Value = "05:21"
$ws.Cells.Item(11, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(11, 4).Value = 63
$ws.Cells.Item(12, 2).Value = "05:22"
$ws.Cells.Item(12, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(12, 4).Value = 90
$ws.Cells.Item(13, 1).Value = "04:18:53"
$ws.Cells.Item(13, 2).Value = "05:34"
$ws.Cells.Item(13, 3).Value = "215B_EL PATO"
$ws.Cells.Item(13, 4).Value = 76
$ws.Cells.Item(13, 5).Value = "LP1912"
$ws.Cells.Item(14, 1).Value = "04:18:53"
$ws.Cells.Item(14, 2).Value = "05:46"
$ws.Cells.Item(14, 3).Value = "15_ABASTO"
$ws.Cells.Item(14, 4).Value = 88
$ws.Cells.Item(14, 5).Value = "LP1912"
$ws.Cells.Item(15, 1).Value = "04:18:53"
$ws.Cells.Item(15, 2).Value = "05:53"
$ws.Cells.Item(15, 3).Value = "10_OLMOS"
$ws.Cells.Item(15, 4).Value = 95
$ws.Cells.Item(15, 5).Value = "LP1912"
$ws.Cells.Item(16, 1).Value = "04:18:53"
$ws.Cells.Item(16, 2).Value = "06:05"
$ws.Cells.Item(16, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(16, 4).Value = 107
$ws.Cells.Item(16, 5).Value = "LP1912"
$ws.Cells.Item(17, 1).Value = "04:18:53"
$ws.Cells.Item(17, 2).Value = "06:11"
$ws.Cells.Item(17, 3).Value = "215A_EL PATO"
$ws.Cells.Item(17, 4).Value = 113
$ws.Cells.Item(17, 5).Value = "LP1912"
$ws.Cells.Item(18, 1).Value = "04:18:53"
$ws.Cells.Item(18, 2).Value = "06:13"
$ws.Cells.Item(18, 3).Value = "225_HARAS DEL SUR"
$ws.Cells.Item(18, 4).Value = 115
$ws.Cells.Item(18, 5).Value = "LP1912"

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 1).Value = "Última actualización: 04:18:53"
$ws.Cells.Item(3, 1).Value = "Total filas: 4"
$ws.Cells.Item(6, 1).Value = "04:18:53"
$ws.Cells.Item(6, 2).Value = "04:45"
$ws.Cells.Item(6, 4).Value = 27
$ws.Cells.Item(7, 2).Value = "04:46"
$ws.Cells.Item(7, 3).Value = "215A_EL PATO"
$ws.Cells.Item(7, 4).Value = 54
$ws.Cells.Item(8, 1).Value = "04:18:53"
$ws.Cells.Item(8, 2).Value = "05:34"
$ws.Cells.Item(8, 3).Value = "215B_EL PATO"
$ws.Cells.Item(8, 4).Value = 76
$ws.Cells.Item(8, 5).Value = "LP1912"
$ws.Cells.Item(9, 1).Value = "04:18:53"
$ws.Cells.Item(9, 2).Value = "06:11"
$ws.Cells.Item(9, 3).Value = "215A_EL PATO"
$ws.Cells.Item(9, 4).Value = 113
$ws.Cells.Item(9, 5).Value = "LP1912"

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 1).Value = "Última actualización: 04:18:53"
$ws.Cells.Item(3, 1).Value = "Total filas: 3"
$ws.Cells.Item(6, 1).Value = "04:18:53"
$ws.Cells.Item(6, 2).Value = "05:43"
$ws.Cells.Item(6, 4).Value = 85
$ws.Cells.Item(7, 1).Value = "03:52:29"
$ws.Cells.Item(7, 2).Value = "05:44"
$ws.Cells.Item(7, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(7, 4).Value = 112
$ws.Cells.Item(7, 5).Value = "L6173"
$ws.Cells.Item(8, 1).Value = "04:18:53"
$ws.Cells.Item(8, 2).Value = "06:08"
$ws.Cells.Item(8, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(8, 4).Value = 110
$ws.Cells.Item(8, 5).Value = "L6173"
